$d = $word.ActiveDocument

# 1. Update the "Tipe soal yang didukung" instruction line to mention the
#    new NUMERICAL_INPUT question type.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "4. Tipe soal yang didukung:*") {
        $p.Range.Text = "4. Tipe soal yang didukung: MULTIPLE_CHOICE, MULTIPLE_SELECTION, TRUE_FALSE, MATCHING, ORDERING, ESSAY, NUMERICAL_INPUT"
    }
}

# 2. Add a new example row for NUMERICAL_INPUT right after the ESSAY row
#    (and before the trailing blank template rows) in the question table.
$t = $d.Tables.Item(1)

$essayRow = $null
foreach ($row in $t.Rows) {
    if ($row.Cells.Item(1).Range.Text -like "ESSAY*") {
        $essayRow = $row
    }
}

$beforeRow = $t.Rows.Item($essayRow.Index + 1)
$newRow = $t.Rows.Add($beforeRow)
$idx = $newRow.Index

$t.Cell($idx, 1).Range.Text = "NUMERICAL_INPUT"

$questionCell = $t.Cell($idx, 2)
$questionCell.Range.Text = "Hitung nilai dari "
$qr = $questionCell.Range
$qr.Collapse(0)
$qr.InsertAfter("`r`$\frac{3}{4} + \frac{2}{5}`$")

$t.Cell($idx, 3).Range.Text = "-"
$t.Cell($idx, 4).Range.Text = "1.15"
$t.Cell($idx, 5).Range.Text = "15"
